$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 23

$ws.Cells.Item($row, 1).Value = 4
$ws.Cells.Item($row, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item($row, 3).Value = "Los Lagos"

# Date value (serial 44890 -> 2022-11-25), using the same style as the other date cells (D column)
$ws.Cells.Item($row, 4).Value = 44890
$ws.Cells.Item($row, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Cells.Item($row, 5).Value = 10
$ws.Cells.Item($row, 6).Value = "Fruta"
$ws.Cells.Item($row, 7).Value = 100101
$ws.Cells.Item($row, 8).Value = "Berries"
$ws.Cells.Item($row, 9).Value = 100101001
$ws.Cells.Item($row, 10).Value = "Arándano (blue)"
$ws.Cells.Item($row, 11).Value = "Sin especificar"
$ws.Cells.Item($row, 12).Value = "Primera"
$ws.Cells.Item($row, 13).Value = 300
$ws.Cells.Item($row, 14).Value = 8500
$ws.Cells.Item($row, 15).Value = 9000
$ws.Cells.Item($row, 16).Value = 8750
$ws.Cells.Item($row, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item($row, 18).Value = "Provincia de Curicó"
$ws.Cells.Item($row, 19).Value = 4375
$ws.Cells.Item($row, 20).Value = 2
